$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "57.440.07"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -3.82%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.921.48"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  -0.09%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "547.58"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.68%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "130.18"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.91%  "

$ws.Range("E7").Value = "  -0.25%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.513"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +3.11%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.915.65"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("E10").Value = "  -2.03%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "4.74"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -6.77%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.445"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.30%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000221"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.01%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "32.80"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.57%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.121"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.62%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.402.31"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "6.84"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +7.63%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.914.86"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -2.69%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "57.463.00"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -4.01%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "415.71"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.10"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.88%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.682"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +3.19%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.93"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.03"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "79.56"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("E26").Value = "  -0.01%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.38%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.15%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.44"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.42%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.98"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "25.14"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.95"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.12%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0961"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.77%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.49%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.944"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.21%  "

$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("E37").Value = "  -3.66%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "8.68"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.92%  "

$ws.Range("E39").Value = "  +3.62%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.54"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.04%  "

$ws.Range("E41").Value = "  -0.99%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0344"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "370.98"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.90%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.667.60"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.46%  "

$ws.Range("E45").Value = "  -0.04%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "123.00"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "

$ws.Range("E47").Value = "  +2.31%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("E49").Value = "  -0.55%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "23.10"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.00%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.00"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "
